$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Swap the table style on the table on slide 6 to the new style GUID.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{F09C5FC1-FDB3-49F6-B012-7160253C8489}")
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the deck's theme colours from the "Integral" palette to the
#    "Office Theme" palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
# ---------------------------------------------------------------------------
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
